# Add a new worksheet "2025-06-06" after "2025-06-04", with the same
# header row layout used by the existing sheets.

$wb = $excel.ActiveWorkbook

# Locate the "2025-06-04" sheet so the new sheet can be placed right after it.
$afterSheet = $wb.Worksheets.Item("2025-06-04")

$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$newSheet.Name = "2025-06-06"

$headers = @("FECHA", "TEAM", "AGENTE", "NÚMERO", "SERVICIO", "PUNTOS", "CUENTA", "DIRECCIÓN", "ZIP CODE")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $newSheet.Cells.Item(1, $i + 1).Value = $headers[$i]
}
